# Knight Unit Sprite.xlsx - add a "Dimensions" column (F) with "96X96" for
# every sprite-row (3-26), matching the header/body styling already used by
# the rest of the table, then update the saved selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- New header cell F2: "Dimensions" (bold, same look as A2:E2) ---------
$ws.Range("F2").Value = "Dimensions"
$ws.Range("F2").Font.Bold = $true

# --- New body cells F3:F26: "96X96", right aligned (same look as J7/J9/J11)
$ws.Range("F3:F26").Value = "96X96"
$ws.Range("F3:F26").HorizontalAlignment = -4152   # xlRight

# --- Column F width, sized to fit the new content -------------------------
$ws.Columns.Item(6).ColumnWidth = 10.6

# --- Update the sheet's saved selection/view -------------------------------
$ws.Range("F3:F26").Select() | Out-Null
